$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 47, pushing existing rows 47:104 down to 48:105.
$ws.Rows.Item(47).Insert()

# Populate the new row 47 with the new data record.
$ws.Cells.Item(47, 1).Value = 4
$ws.Cells.Item(47, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(47, 3).Value = "Los Lagos"
$ws.Cells.Item(47, 4).Value = 44413
$ws.Cells.Item(47, 5).Value = 10
$ws.Cells.Item(47, 6).Value = 100112032
$ws.Cells.Item(47, 7).Value = "Zapallo italiano"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 200
$ws.Cells.Item(47, 11).Value = 11000
$ws.Cells.Item(47, 12).Value = 12000
$ws.Cells.Item(47, 13).Value = 11500
$ws.Cells.Item(47, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(47, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(47, 16).Value = 230
$ws.Cells.Item(47, 17).Value = 50
$ws.Cells.Item(47, 18).Value = "Hortaliza"
